# Auto-update draw results: append the 2025-12-22 "Pick 4" draw as a new
# row at the bottom of the Results sheet, mirroring every prior row's
# layout: Date | Game | Phase | Result | InsertedAt.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first free row below the existing data (classic "Ctrl+Up from
# the bottom of the sheet" idiom), so the new draw lands right after the
# last one regardless of exactly how many rows already exist.
$lastRow = $ws.Range("A1048576").End(-4162).Row
$row = $lastRow + 1

# A (Date) and C (Phase) look like a date / a plain integer respectively,
# so a bare assignment would get auto-converted to a date serial number /
# numeric value. A leading apostrophe forces Excel to store the literal
# text instead, matching every other row in the column. Resetting the
# cell style back to "Normal" afterwards clears the "quote prefix" text
# flag that gets applied when text is forced this way, so the cell keeps
# the sheet's plain, unstyled look (same as all the other data cells).
$ws.Range("A$row").Value = "'2025-12-22"
$ws.Range("A$row").Style = "Normal"

$ws.Range("B$row").Value = "Pick 4"

$ws.Range("C$row").Value = "'251222"
$ws.Range("C$row").Style = "Normal"

$ws.Range("D$row").Value = "5-7-9-9"

$ws.Range("E$row").Value = "2025-12-22T21:41:00.267+04:00"
